# Update "想去人数" (want-to-go count) values in F column on sheets "展览" and "全部类型"
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 449
$ws1.Range("F3").Value = 28
$ws1.Range("F4").Value = 37
$ws1.Range("F5").Value = 4939
$ws1.Range("F6").Value = 169
$ws1.Range("F7").Value = 76
$ws1.Range("F8").Value = 273

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 449
$ws4.Range("F7").Value = 28
$ws4.Range("F8").Value = 37
$ws4.Range("F9").Value = 4939
$ws4.Range("F10").Value = 169
$ws4.Range("F11").Value = 76
$ws4.Range("F13").Value = 273
